$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.9999999992646067,
    0.99547375265234017,
    0.97494050540500299,
    0.96683286039723559,
    0.95918810983866121,
    0.94099166954989433,
    0.93735635355740898,
    0.93286898584076616,
    0.92901533009105153,
    0.9181839202086477,
    0.91658741361611784,
    0.91380706479607543,
    0.90251818920744187,
    0.89835032600394737,
    0.89575855250618508,
    0.89325187128124872,
    0.88954371055450232,
    0.88843475668407956,
    0.994420456022445,
    0.98730313823370364,
    0.98590460135644697,
    0.98464008707227502,
    0.96631096497866331,
    0.94998018602123535,
    0.94352298582423444,
    0.93933374745755294,
    0.93727329781983282,
    0.92870870588505072,
    0.92718205875702375,
    0.92768779893955244,
    0.93751256273917338,
    0.94187141228328353,
    0.94867369977208238
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

